$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.320.73"
$ws.Range("E2").Value = "  -0.82%  "
$ws.Range("D3").Value = "1.651.06"
$ws.Range("E3").Value = "  -0.55%  "
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").Value = "212.92"
$ws.Range("E5").Value = "  -0.78%  "
$ws.Range("D6").Value = "0.510"
$ws.Range("E6").Value = "  -0.50%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("D8").Value = "23.32"
$ws.Range("E8").Value = "  +0.17%  "
$ws.Range("D10").Value = "0.0613"
$ws.Range("E10").Value = "  -1.01%  "
$ws.Range("D11").Value = "0.0895"
$ws.Range("E11").Value = "  +2.30%  "
$ws.Range("D12").Value = "1.884.97"
$ws.Range("E12").Value = "  -0.52%  "
$ws.Range("D13").Value = "1.649.76"
$ws.Range("E13").Value = "  -0.61%  "
$ws.Range("E14").Value = "  -1.24%  "
$ws.Range("D15").Value = "0.566"
$ws.Range("E15").Value = "  +3.20%  "
$ws.Range("D16").Value = "65.41"
$ws.Range("E16").Value = "  -0.46%  "
$ws.Range("D17").Value = "27.316.20"
$ws.Range("E17").Value = "  -0.82%  "
$ws.Range("D18").Value = "229.87"
$ws.Range("E18").Value = "  -6.42%  "
$ws.Range("D19").Value = "0.0₃0726"
$ws.Range("E19").Value = "  -0.53%  "
$ws.Range("D20").Value = "7.38"
$ws.Range("E20").Value = "  -1.11%  "
$ws.Range("E21").Value = "  -0.04%  "
$ws.Range("D22").Value = "4.35"
$ws.Range("E22").Value = "  -2.71%  "
$ws.Range("D23").Value = "9.43"
$ws.Range("E23").Value = "  +1.23%  "
$ws.Range("D24").Value = "2.04"
$ws.Range("E24").Value = "  +0.89%  "
$ws.Range("D25").Value = "147.21"
$ws.Range("E25").Value = "  +0.82%  "
$ws.Range("D26").Value = "7.05"
$ws.Range("E26").Value = "  -1.56%  "
$ws.Range("D27").Value = "15.78"
$ws.Range("E27").Value = "  -2.60%  "
$ws.Range("E28").Value = "  -0.14%  "
$ws.Range("D29").Value = "0.111"
$ws.Range("E29").Value = "  -0.03%  "
$ws.Range("E30").Value = "  -0.38%  "
$ws.Range("D31").Value = "1.18"
$ws.Range("E31").Value = "  -4.40%  "
$ws.Range("E32").Value = "  -1.18%  "
$ws.Range("D33").Value = "1.427.55"
$ws.Range("E33").Value = "  -0.87%  "
$ws.Range("E34").Value = "  -0.11%  "
$ws.Range("E35").Value = "  +0.89%  "
$ws.Range("E36").Value = "  -0.33%  "
$ws.Range("D37").Value = "0.903"
$ws.Range("E37").Value = "  -2.68%  "
$ws.Range("D38").Value = "0.569"
$ws.Range("E38").Value = "  -1.68%  "
$ws.Range("E39").Value = "  -0.44%  "
$ws.Range("E40").Value = "  +0.74%  "
$ws.Range("E41").Value = "  -0.05%  "
$ws.Range("D42").Value = "5.55"
$ws.Range("E42").Value = "  +2.62%  "
$ws.Range("D43").Value = "64.85"
$ws.Range("E43").Value = "  -6.04%  "
$ws.Range("E44").Value = "  +0.62%  "
$ws.Range("D45").Value = "0.786"
$ws.Range("E45").Value = "  -0.86%  "
$ws.Range("D46").Value = "1.793.50"
$ws.Range("E46").Value = "  -0.54%  "
$ws.Range("E47").Value = "  -2.44%  "
$ws.Range("D48").Value = "87.84"
$ws.Range("E48").Value = "  -0.85%  "
$ws.Range("D49").Value = "0.0₆0106"
$ws.Range("E49").Value = "  -1.73%  "
$ws.Range("E50").Value = "  -0.04%  "
$ws.Range("D51").Value = "7.70"
$ws.Range("E51").Value = "  -1.45%  "
